$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A (new "data" and "loja" columns),
# shifting the existing nome/modelo/preco/politica/full/tipo/link columns
# from A:G to C:I.
$ws.Range("A1:B1").EntireColumn.Insert()

# New header labels
$ws.Range("A1").Value = "data"
$ws.Range("B1").Value = "loja"

# Match the bold/centered header style used by the other header cells
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new "data" and "loja" values for every data row
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = "30/07/2024"
    $ws.Cells.Item($r, 2).Value = "tudo.som"
}

# Update the tracking_id portion of every link (now in column I)
$ws.Range("I2:I9").Replace("61aad6a6-bae6-4208-9ec0-1a7abf86053b", "aa68dec7-8970-4a14-b0a1-3f9cb98e5acd")
